$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Right" total (B11: Marking row's Right value) and "Total" row values
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 130
$ws.Range("E12").Value = "130/140"
